$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New MarketBeat snapshot (10th run): shift history right by inserting
# three new date columns in front of the old "B" column (old B stays put and
# becomes the new B1 value; old C,D,E shift to F,G,H). ---
$ws.Columns("C:E").Insert()

# Header row: old B1 ("Jun_17") needs to move into the freshly vacated E1
# slot before we overwrite B1 with the newest date.
$ws.Range("E1").Value = "Jun_17"
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# The three newly inserted columns have no data yet for the existing
# analyst rows (2-27) -- mark them "UN" (unchanged), matching every other
# untouched rating cell in the sheet.
$ws.Range("C2:E27").Value = "UN"

# --- Added new group: two more analyst rows at the bottom of the table ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
